$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Insert "Professor " before "Francisco Sepulveda" (Supervisor line)
# ---------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("Francisco Sepulveda", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $find.Start

# Duplicate the first character of "Francisco Sepulveda" (carries the exact
# Arial/23/23 run formatting we need) and insert the copy right before the
# original text, leaving the original untouched.
$srcChar = $d.Range($startPos, $startPos + 1)
$ft = $srcChar.FormattedText
$insPoint = $d.Range($startPos, $startPos)
$insPoint.FormattedText = $ft

# Retext the inserted copy to be "Professor " -- this merges it into the
# following run (same formatting), so force a clean re-split by toggling a
# character property on just the inserted span.
$newCopy = $d.Range($startPos, $startPos + 1)
$newCopy.Text = "Professor "
$profRange = $d.Range($startPos, $startPos + 10)
$profRange.Font.Bold = $true
$profRange.Font.Bold = $false

# ---------------------------------------------------------------
# 2 & 3. Move the "_GoBack" bookmark from after "Intro" to the end of the
#        "Abstract" paragraph, and drop the trailing "." run there.
# ---------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$findAbstract = $d.Content
$foundAbstract = $findAbstract.Find.Execute("Abstract", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($findAbstract.End, $findAbstract.End)
$d.Bookmarks.Add("_GoBack", $target)

$findPeriod = $d.Content
$foundPeriod = $findPeriod.Find.Execute("Abstract.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$periodRange = $d.Range($findPeriod.End - 1, $findPeriod.End)
$periodRange.Delete()
